$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "30-2=28"
$t.Cell(1,2).Range.Text = "72-38=34"
$t.Cell(1,3).Range.Text = "61-2=59"
$t.Cell(1,4).Range.Text = "21-7=14"
$t.Cell(1,5).Range.Text = "58+6=64"
$t.Cell(2,1).Range.Text = "14+18=32"
$t.Cell(2,2).Range.Text = "29+4=33"
$t.Cell(2,3).Range.Text = "89+6=95"
$t.Cell(2,4).Range.Text = "19+67=86"
$t.Cell(2,5).Range.Text = "77+17=94"
$t.Cell(3,1).Range.Text = "55+36=91"
$t.Cell(3,2).Range.Text = "62+19=81"
$t.Cell(3,3).Range.Text = "55-18=37"
$t.Cell(3,4).Range.Text = "62-14=48"
$t.Cell(3,5).Range.Text = "39+8=47"
$t.Cell(4,1).Range.Text = "7+29=36"
$t.Cell(4,2).Range.Text = "93-25=68"
$t.Cell(4,3).Range.Text = "93-36=57"
$t.Cell(4,4).Range.Text = "72-5=67"
$t.Cell(4,5).Range.Text = "53-19=34"
$t.Cell(5,1).Range.Text = "60-38=22"
$t.Cell(5,2).Range.Text = "2+49=51"
$t.Cell(5,3).Range.Text = "2+59=61"
$t.Cell(5,4).Range.Text = "81-13=68"
$t.Cell(5,5).Range.Text = "7+24=31"
$t.Cell(6,1).Range.Text = "63-57=6"
$t.Cell(6,2).Range.Text = "4+48=52"
$t.Cell(6,3).Range.Text = "81-66=15"
$t.Cell(6,4).Range.Text = "35+17=52"
$t.Cell(6,5).Range.Text = "56+7=63"
$t.Cell(7,1).Range.Text = "17+27=44"
$t.Cell(7,2).Range.Text = "28-9=19"
$t.Cell(7,3).Range.Text = "70-18=52"
$t.Cell(7,4).Range.Text = "60-27=33"
$t.Cell(7,5).Range.Text = "73-68=5"
$t.Cell(8,1).Range.Text = "50-14=36"
$t.Cell(8,2).Range.Text = "25-16=9"
$t.Cell(8,3).Range.Text = "68-19=49"
$t.Cell(8,4).Range.Text = "84-16=68"
$t.Cell(8,5).Range.Text = "57+7=64"
$t.Cell(9,1).Range.Text = "85-47=38"
$t.Cell(9,2).Range.Text = "25+69=94"
$t.Cell(9,3).Range.Text = "6+86=92"
$t.Cell(9,4).Range.Text = "55+16=71"
$t.Cell(9,5).Range.Text = "95-17=78"
$t.Cell(10,1).Range.Text = "8+16=24"
$t.Cell(10,2).Range.Text = "76-69=7"
$t.Cell(10,3).Range.Text = "60-49=11"
$t.Cell(10,4).Range.Text = "46-17=29"
$t.Cell(10,5).Range.Text = "73-48=25"
$t.Cell(11,1).Range.Text = "87-9=78"
$t.Cell(11,2).Range.Text = "71-45=26"
$t.Cell(11,3).Range.Text = "15+6=21"
$t.Cell(11,4).Range.Text = "96-29=67"
$t.Cell(11,5).Range.Text = "21-14=7"
$t.Cell(12,1).Range.Text = "11-4=7"
$t.Cell(12,2).Range.Text = "36+35=71"
$t.Cell(12,3).Range.Text = "92-17=75"
$t.Cell(12,4).Range.Text = "30-23=7"
$t.Cell(12,5).Range.Text = "63-9=54"
$t.Cell(13,1).Range.Text = "13+38=51"
$t.Cell(13,2).Range.Text = "8+77=85"
$t.Cell(13,3).Range.Text = "5+19=24"
$t.Cell(13,4).Range.Text = "60-26=34"
$t.Cell(13,5).Range.Text = "43+48=91"
$t.Cell(14,1).Range.Text = "41-33=8"
$t.Cell(14,2).Range.Text = "82-34=48"
$t.Cell(14,3).Range.Text = "28+47=75"
$t.Cell(14,4).Range.Text = "91-46=45"
$t.Cell(14,5).Range.Text = "83-6=77"
$t.Cell(15,1).Range.Text = "77-29=48"
$t.Cell(15,2).Range.Text = "43-16=27"
$t.Cell(15,3).Range.Text = "33-25=8"
$t.Cell(15,4).Range.Text = "70-69=1"
$t.Cell(15,5).Range.Text = "94-67=27"
$t.Cell(16,1).Range.Text = "68+24=92"
$t.Cell(16,2).Range.Text = "34+58=92"
$t.Cell(16,3).Range.Text = "78+16=94"
$t.Cell(16,4).Range.Text = "81-48=33"
$t.Cell(16,5).Range.Text = "50-24=26"
$t.Cell(17,1).Range.Text = "53-35=18"
$t.Cell(17,2).Range.Text = "47-39=8"
$t.Cell(17,3).Range.Text = "4+29=33"
$t.Cell(17,4).Range.Text = "91-39=52"
$t.Cell(17,5).Range.Text = "96-28=68"
$t.Cell(18,1).Range.Text = "16+25=41"
$t.Cell(18,2).Range.Text = "47+29=76"
$t.Cell(18,3).Range.Text = "64-38=26"
$t.Cell(18,4).Range.Text = "37+28=65"
$t.Cell(18,5).Range.Text = "80-37=43"
$t.Cell(19,1).Range.Text = "76+18=94"
$t.Cell(19,2).Range.Text = "91-85=6"
$t.Cell(19,3).Range.Text = "7+79=86"
$t.Cell(19,4).Range.Text = "62-53=9"
$t.Cell(19,5).Range.Text = "90-82=8"
$t.Cell(20,1).Range.Text = "29+34=63"
$t.Cell(20,2).Range.Text = "90-43=47"
$t.Cell(20,3).Range.Text = "19+33=52"
$t.Cell(20,4).Range.Text = "7+64=71"
$t.Cell(20,5).Range.Text = "70-52=18"